$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$ws.Range("A2").Value = "2025-11-02 01:54:19"
$ws.Range("A3").Value = "2025-11-02 01:54:19"
$ws.Range("A4").Value = "2025-11-02 01:54:19"
$ws.Range("A5").Value = "2025-11-02 01:54:19"
$ws.Range("A6").Value = "2025-11-02 01:54:19"
